$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "'67.556.32"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.19%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "'3.783.24"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.28%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = "'597.17"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.23%  "

# Row 6: Solana
$ws.Cells.Item(6, 4).Value = "'164.60"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.57%  "

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = "  +0.04%  "

# Row 8: XRP
$ws.Cells.Item(8, 4).Value = "'0.514"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.04%  "

# Row 9: Dogecoin
$ws.Cells.Item(9, 5).Value = "  -1.19%  "

# Row 10: Cardano
$ws.Cells.Item(10, 5).Value = "  +0.22%  "

# Row 11: Toncoin
$ws.Cells.Item(11, 4).Value = "'6.40"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.20%  "

# Row 12: ShibaInu
$ws.Cells.Item(12, 4).Value = "'0.0000247"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -2.40%  "

# Row 13: Avalanche
$ws.Cells.Item(13, 4).Value = "'35.56"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.56%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "'4.422.21"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.09%  "

# Row 15: WrappedEther
$ws.Cells.Item(15, 4).Value = "'3.818.86"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.95%  "

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = "'67.636.51"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.01%  "

# Row 17: Chainlink
$ws.Cells.Item(17, 4).Value = "'18.25"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.85%  "

# Row 18: TRON
$ws.Cells.Item(18, 5).Value = "  +1.67%  "

# Row 19: Polkadot
$ws.Cells.Item(19, 4).Value = "'7.01"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.50%  "

# Row 20: BitcoinCash
$ws.Cells.Item(20, 4).Value = "'460.33"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.22%  "

# Row 21: Uniswap
$ws.Cells.Item(21, 4).Value = "'9.74"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.50%  "

# Row 22: Polygon
$ws.Cells.Item(22, 5).Value = "  -0.33%  "

# Row 23: PEPE
$ws.Cells.Item(23, 2).Value = "PEPE"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(23, 4).Value = "'0.0000145"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -6.60%  "

# Row 24: Litecoin
$ws.Cells.Item(24, 2).Value = "Litecoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(24, 4).Value = "'82.36"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.11%  "

# Row 25: InternetComputer(DFINITY)
$ws.Cells.Item(25, 4).Value = "'11.96"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.74%  "

# Row 26: Fetch.AI
$ws.Cells.Item(26, 4).Value = "'2.08"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.60%  "

# Row 27: Dai
$ws.Cells.Item(27, 5).Value = "  -0.01%  "

# Row 28: RenderToken
$ws.Cells.Item(28, 5).Value = "  -0.26%  "

# Row 29: WrappedeETH
$ws.Cells.Item(29, 4).Value = "'3.933.44"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.16%  "

# Row 30: NEARProtocol
$ws.Cells.Item(30, 5).Value = "  +1.84%  "

# Row 31: PancakeSwap
$ws.Cells.Item(31, 5).Value = "  -4.48%  "

# Row 32: ImmutableX
$ws.Cells.Item(32, 4).Value = "'2.18"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.31%  "

# Row 33: EthereumClassic
$ws.Cells.Item(33, 4).Value = "'28.88"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -2.62%  "

# Row 35: Aptos
$ws.Cells.Item(35, 4).Value = "'8.96"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.42%  "

# Row 36: Hedera
$ws.Cells.Item(36, 4).Value = "'0.0985"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.52%  "

# Row 37: Kaspa
$ws.Cells.Item(37, 5).Value = "  +0.25%  "

# Row 38: Mantle
$ws.Cells.Item(38, 4).Value = "'0.987"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.44%  "

# Row 39: dogwifhat
$ws.Cells.Item(39, 2).Value = "dogwifhat"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(39, 4).Value = "'3.20"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -5.12%  "

# Row 40: Filecoin
$ws.Cells.Item(40, 2).Value = "Filecoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(40, 4).Value = "'5.72"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.88%  "

# Row 41: FirstDigitalUSD
$ws.Cells.Item(41, 5).Value = "  +0.06%  "

# Row 43: Arweave
$ws.Cells.Item(43, 4).Value = "'43.83"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.97%  "

# Row 44: OKB
$ws.Cells.Item(44, 4).Value = "'47.49"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.13%  "

# Row 45: TheGraph
$ws.Cells.Item(45, 4).Value = "'0.296"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.91%  "

# Row 46: Monero
$ws.Cells.Item(46, 4).Value = "'151.58"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.53%  "

# Row 47: Cosmos
$ws.Cells.Item(47, 4).Value = "'8.30"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.19%  "

# Row 48: ONDO
$ws.Cells.Item(48, 5).Value = "  +7.22%  "

# Row 49: EnergySwap
$ws.Cells.Item(49, 4).Value = "'27.14"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.81%  "

# Row 50: Bittensor
$ws.Cells.Item(50, 4).Value = "'393.70"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.04%  "

# Row 51: Stacks
$ws.Cells.Item(51, 4).Value = "'1.84"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.35%  "

